# Update DISH yearly financials: add FY2018 (period ending 2018-12-31) as
# a new first data column (column D), shifting the existing years right
# by one column, and correct a handful of values that changed between
# the two extracts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert a new column before column D. This shifts D:K -> E:L,
#    carrying formatting/number-formats/styles along with the cells.
# ---------------------------------------------------------------------
$ws.Range("D1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. Income statement (rows 7-35) - new column D values for FY2018
# ---------------------------------------------------------------------
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 13621300
$ws.Range("D9").Value = 9609800
$ws.Range("D10").Value = 4011500
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -3700
$ws.Range("D15").Value = 712000
$ws.Range("D17").Value = 11469900
$ws.Range("D18").Value = 2151400
$ws.Range("D20").Value = 52800
$ws.Range("D21").Value = 2916200
$ws.Range("D22").Value = 15000
$ws.Range("D23").Value = 2189200
$ws.Range("D24").Value = 533700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1655500
$ws.Range("D27").Value = 1575100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -52800
$ws.Range("D33").Value = 1575100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1575100

# ---------------------------------------------------------------------
# 3. Balance sheet (rows 38-77)
# ---------------------------------------------------------------------
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 887300
$ws.Range("D42").Value = 1181500
$ws.Range("D43").Value = 639900
$ws.Range("D44").Value = 290700
$ws.Range("D45").Value = 289800
$ws.Range("D46").Value = 3289200
$ws.Range("D47").Value = 119000
$ws.Range("D48").Value = 1928200
$ws.Range("D49").Value = 24754400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 496200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 30587000
$ws.Range("D57").Value = 233800
$ws.Range("D58").Value = 1342000
$ws.Range("D59").Value = 3200400
$ws.Range("D60").Value = 4776100
$ws.Range("D61").Value = 13810800
$ws.Range("D62").Value = 2945800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 21991300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 5212800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 8595700
$ws.Range("D77").Value = 0

# ---------------------------------------------------------------------
# 4. Cash flow statement (rows 80-102)
# ---------------------------------------------------------------------
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 1575100
$ws.Range("D83").Value = 712000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 2517800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -1134500
$ws.Range("D101").Value = "NA"

# Rows 91, 94, 102 were not a pure shift: several existing values were
# also corrected. Set the full D:K range for these rows explicitly.
$ws.Range("D91:K91").Value = @(-393900, -431800, -614100, -735000, -1001900, -1253500, -945300, -760200)
$ws.Range("D94:K94").Value = @(-1975300, -6521600, -1737100, -8062100, -963100, -3021100, -3019200, -2695300)
$ws.Range("D102:K102").Value = @(-592000, -3845300, 4271100, -6051200, 2395300, 1126300, 3003600, -70500)

# ---------------------------------------------------------------------
# 5. New column D needs the same number formats as the columns it was
#    copied from (dates -> style 2, numbers -> style 3).
# ---------------------------------------------------------------------
$ws.Range("D7,D38,D80").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D8:D10,D12:D35,D41:D54,D57:D77,D81,D83:D102").NumberFormat = "#,##0"

$wb.Save()
